$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "testExplicitLogin"
$ws.Range("B3").Value = "PASS"
$ws.Range("C3").Value = "16_05_2018_12_36_00"
$ws.Range("D3").Value = "CHROME"

$ws.Range("A4").Value = "testExplicitLogin"
$ws.Range("B4").Value = "FAIL"
$ws.Range("C4").Value = "16_05_2018_12_38_01"
$ws.Range("D4").Value = "CHROME"
